# Apply the captured commit:
#  1. Re-style the three data tables (slides 14-16) to use the
#     "Medium Style 2 - Accent 1" built-in table style instead of the
#     custom "Table_0" style that used to be the deck's default.
#  2. Swap the presentation's applied theme palette from the pink/violet
#     "Integral" scheme back to the default "Office" scheme (this is the
#     look a user gets after picking the plain "Office Theme" design).

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$oldStyleId = "{8A9E8664-533E-48BE-9641-90AA031D526C}"
$newStyleId = "{B9C3F706-A7D2-4533-B767-38784AC432F5}"

foreach ($slideIdx in 14,15,16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Theme colors ---------------------------------------------------
# Office theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as VBA RGB() values, applied via the presentation's live
# ThemeColorScheme (backs the deck's single shared theme part).
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
